$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2999.2856
$ws.Range("I76").Value = 2999.2307
$ws.Range("K76").Value = 2999.2307
$ws.Range("M76").Value = -2684.2307
$ws.Range("H79").Value = 2999.2856
$ws.Range("I79").Value = 2999.2307
$ws.Range("K79").Value = 2999.2307
$ws.Range("M79").Value = -1907.2307
$ws.Range("H131").Value = 333337500
$ws.Range("I131").Value = 500005000
$ws.Range("J131").Value = 2500
$ws.Range("K131").Value = 1500015000
$ws.Range("L131").Value = 7500
$ws.Range("M131").Value = -1500009960
$ws.Range("N131").Value = -17580
$ws.Range("H137").Value = 1181.6765
$ws.Range("I137").Value = 945.069
$ws.Range("J137").Value = 2554
$ws.Range("K137").Value = 2835.207
$ws.Range("L137").Value = 7662
$ws.Range("M137").Value = -285.2069999999999
$ws.Range("N137").Value = -12762
$ws.Range("H138").Value = 2327913.5
$ws.Range("I138").Value = 1072.9678
$ws.Range("J138").Value = 3639405.5
$ws.Range("K138").Value = 3218.9034
$ws.Range("L138").Value = 10918216.5
$ws.Range("M138").Value = 1921.0966
$ws.Range("N138").Value = -10928496.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4781.744
$ws.Range("I32").Value = 3204.1216
$ws.Range("J32").Value = 19374.75
$ws.Range("K32").Value = 3204.1216
$ws.Range("L32").Value = 19374.75
$ws.Range("M32").Value = -2917.1216
$ws.Range("N32").Value = -19948.75
$ws.Range("H45").Value = 1994.5625
$ws.Range("I45").Value = 1514.2858
$ws.Range("J45").Value = 2368.111
$ws.Range("K45").Value = 1514.2858
$ws.Range("L45").Value = 2368.111
$ws.Range("M45").Value = -1137.2858
$ws.Range("N45").Value = -3122.111
$ws.Range("H52").Value = 19666.666
$ws.Range("J52").Value = 19666.666
$ws.Range("L52").Value = 19666.666
$ws.Range("N52").Value = -20302.666
$ws.Range("H61").Value = 1335.9231
$ws.Range("I61").Value = 1207.5135
$ws.Range("J61").Value = 1652.6666
$ws.Range("K61").Value = 1207.5135
$ws.Range("L61").Value = 1652.6666
$ws.Range("M61").Value = -995.5135
$ws.Range("N61").Value = -2076.6666
$ws.Range("H74").Value = 23016.957
$ws.Range("I74").Value = 32115.312
$ws.Range("J74").Value = 2220.7144
$ws.Range("K74").Value = 32115.312
$ws.Range("L74").Value = 2220.7144
$ws.Range("M74").Value = -31241.312
$ws.Range("N74").Value = -3968.7144
$ws.Range("H77").Value = 23016.957
$ws.Range("I77").Value = 32115.312
$ws.Range("J77").Value = 2220.7144
$ws.Range("K77").Value = 160576.56
$ws.Range("L77").Value = 11103.572
$ws.Range("M77").Value = -156208.56
$ws.Range("N77").Value = -19839.572
$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996
$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984
$ws.Range("H122").Value = 2078.2666
$ws.Range("I122").Value = 1967.7142
$ws.Range("K122").Value = 5903.142599999999
$ws.Range("M122").Value = -3453.142599999999
$ws.Range("H132").Value = 1402.4154
$ws.Range("I132").Value = 1421.2642
$ws.Range("J132").Value = 1319.1666
$ws.Range("K132").Value = 4263.792600000001
$ws.Range("L132").Value = 3957.4998
$ws.Range("M132").Value = -1733.792600000001
$ws.Range("N132").Value = -9017.4998
$ws.Range("H136").Value = 1335.9231
$ws.Range("I136").Value = 1207.5135
$ws.Range("J136").Value = 1652.6666
$ws.Range("K136").Value = 3622.5405
$ws.Range("L136").Value = 4957.9998
$ws.Range("M136").Value = -1072.5405
$ws.Range("N136").Value = -10057.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 410365.6
$ws.Range("I134").Value = 542423.9
$ws.Range("J134").Value = 3185.7917
$ws.Range("K134").Value = 1627271.7
$ws.Range("L134").Value = 9557.375100000001
$ws.Range("M134").Value = -1624736.7
$ws.Range("N134").Value = -14627.3751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1420.2131
$ws.Range("I31").Value = 882.97437
$ws.Range("J31").Value = 2372.5908
$ws.Range("K31").Value = 882.97437
$ws.Range("L31").Value = 2372.5908
$ws.Range("M31").Value = -587.97437
$ws.Range("N31").Value = -2962.5908
$ws.Range("H34").Value = 1420.2131
$ws.Range("I34").Value = 882.97437
$ws.Range("J34").Value = 2372.5908
$ws.Range("K34").Value = 882.97437
$ws.Range("L34").Value = 2372.5908
$ws.Range("M34").Value = -680.97437
$ws.Range("N34").Value = -2776.5908
$ws.Range("H132").Value = 700395.7
$ws.Range("I132").Value = 1589.3556
$ws.Range("J132").Value = 4631181
$ws.Range("K132").Value = 4768.066800000001
$ws.Range("L132").Value = 13893543
$ws.Range("M132").Value = -2238.066800000001
$ws.Range("N132").Value = -13898603
$ws.Range("H134").Value = 1863.2549
$ws.Range("I134").Value = 1858.119
$ws.Range("J134").Value = 1887.2222
$ws.Range("K134").Value = 5574.357
$ws.Range("L134").Value = 5661.6666
$ws.Range("M134").Value = -3039.357
$ws.Range("N134").Value = -10731.6666
$ws.Range("H141").Value = 68440
$ws.Range("J141").Value = 73889.14
$ws.Range("L141").Value = 73889.14
$ws.Range("N141").Value = -84249.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3321.6667
$ws.Range("I3").Value = 1982.5
$ws.Range("K3").Value = 5947.5
$ws.Range("M3").Value = -5835.5
$ws.Range("H107").Value = 625408.5
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 667089.0600000001
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 2001267.18
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -2005107.18
$ws.Range("H109").Value = 2006.238
$ws.Range("I109").Value = 839.4545000000001
$ws.Range("J109").Value = 3289.7
$ws.Range("K109").Value = 2518.3635
$ws.Range("L109").Value = 9869.099999999999
$ws.Range("M109").Value = -1478.3635
$ws.Range("N109").Value = -11949.1
$ws.Range("H115").Value = 2588.3
$ws.Range("I115").Value = 1027.6666
$ws.Range("J115").Value = 3257.1428
$ws.Range("K115").Value = 3082.9998
$ws.Range("L115").Value = 9771.428400000001
$ws.Range("M115").Value = -1907.9998
$ws.Range("N115").Value = -12121.4284
$ws.Range("H131").Value = 918.03
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 918.03
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2754.09
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12834.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1328.64
$ws.Range("I82").Value = 1542.909
$ws.Range("K82").Value = 1542.909
$ws.Range("M82").Value = -1181.909
$ws.Range("H85").Value = 1328.64
$ws.Range("I85").Value = 1542.909
$ws.Range("K85").Value = 1542.909
$ws.Range("M85").Value = -294.9090000000001
$ws.Range("H132").Value = 4422.7407
$ws.Range("I132").Value = 4362.077
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 13086.231
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -10556.231
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3108.6
$ws.Range("I132").Value = 3787
$ws.Range("J132").Value = 1808.3334
$ws.Range("K132").Value = 11361
$ws.Range("L132").Value = 5425.0002
$ws.Range("M132").Value = -8831
$ws.Range("N132").Value = -10485.0002
